$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: clear stray duplicate value in C2 (row 2 leading y_1 entry)
$ws.Range("C2").ClearContents()

# Corrected forecast values (tiny floating point fix from naive forecaster bug)
$ws.Range("E2").Value = -0.7976031983999876
$ws.Range("E3").Value = -8.396348489509165
$ws.Range("C6").Value = 1.15368307467123
$ws.Range("E6").Value = 3.648892256099967
$ws.Range("C7").Value = 0.2186142574756245
$ws.Range("E7").Value = 0.4006004000999486
$ws.Range("E8").Value = -3.55169094390001
$ws.Range("C9").Value = -0.2262139320475476
$ws.Range("E10").Value = -0.3994003999000184
$ws.Range("C11").Value = 0.3239252862367259
$ws.Range("E12").Value = -0.3994003998999962
$ws.Range("C13").Value = -0.4781004700720182
$ws.Range("E13").Value = 0.8024032016000104
$ws.Range("C14").Value = -1.197849743493751
$ws.Range("E14").Value = -3.161804390400014
$ws.Range("C15").Value = 1.098150690304189
$ws.Range("E15").Value = -2.378486270399993
$ws.Range("C16").Value = -2.376072963557374
$ws.Range("C18").Value = 1.052599339874583
$ws.Range("E18").Value = 1.205410808099971
$ws.Range("C19").Value = -0.7109608111999011
